$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17354691028595
$ws.Range("B1").Value = 2.437618732452393
$ws.Range("D1").Value = 2.365114450454712
$ws.Range("E1").Value = 1.23620343208313
